$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.57"
$ws.Range("E2").Value = "'0.66%"
$ws.Range("D3").Value = "'35.65"
$ws.Range("E3").Value = "'11.23%"
$ws.Range("D4").Value = "'5.089"
$ws.Range("E4").Value = "'1.32%"
$ws.Range("D5").Value = "'0.07787"
$ws.Range("E5").Value = "'-0.24%"
$ws.Range("D6").Value = "'2.263"
$ws.Range("E6").Value = "'-2.31%"
$ws.Range("D7").Value = "'8.077"
$ws.Range("E7").Value = "'1.27%"
$ws.Range("D8").Value = "'4.051"
$ws.Range("E8").Value = "'4.14%"
$ws.Range("D9").Value = "'0.9289"
$ws.Range("E9").Value = "'-0.35%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09277"
$ws.Range("E10").Value = "'-8.67%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1834"
$ws.Range("E11").Value = "'3.49%"
$ws.Range("D12").Value = "'0.08526"
$ws.Range("E12").Value = "'0.73%"
$ws.Range("D13").Value = "'0.03749"
$ws.Range("E13").Value = "'12.21%"
$ws.Range("D14").Value = "'0.09937"
$ws.Range("E14").Value = "'0.41%"
$ws.Range("D15").Value = "'0.001479"
$ws.Range("E15").Value = "'0.25%"
$ws.Range("D16").Value = "'0.005732"
$ws.Range("E16").Value = "'-0.59%"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'2.184"
$ws.Range("E18").Value = "'-0.30%"
$ws.Range("E19").Value = "'2.99%"
$ws.Range("E20").Value = "'-1.35%"
$ws.Range("D21").Value = "'4.610"
$ws.Range("E21").Value = "'7.04%"
$ws.Range("D22").Value = "'0.2240"
$ws.Range("E22").Value = "'7.61%"
$ws.Range("D23").Value = "'0.04676"
$ws.Range("E23").Value = "'0.92%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'1.45%"
$ws.Range("D25").Value = "'0.004532"
$ws.Range("E25").Value = "'3.26%"
$ws.Range("D26").Value = "'0.0001305"
$ws.Range("E26").Value = "'0.91%"
$ws.Range("E27").Value = "'-19.99%"
$ws.Range("D39").Value = "'0.01773"
$ws.Range("E39").Value = "'3.90%"
$ws.Range("E40").Value = "'-0.59%"
$ws.Range("D41").Value = "'0.008000"
$ws.Range("E41").Value = "'3.61%"
$ws.Range("D42").Value = "'0.1415"
$ws.Range("E42").Value = "'0.77%"
$ws.Range("D43").Value = "'0.007889"
$ws.Range("E43").Value = "'-19.28%"
$ws.Range("D44").Value = "'0.002229"
$ws.Range("E44").Value = "'7.45%"
$ws.Range("D45").Value = "'0.009633"
$ws.Range("E45").Value = "'-0.59%"
$ws.Range("D46").Value = "'0.00006199"
$ws.Range("E46").Value = "'1.33%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.93%"
$ws.Range("D48").Value = "'5.278"
$ws.Range("E48").Value = "'106.90%"
$ws.Range("D49").Value = "'0.002699"
$ws.Range("E49").Value = "'35.74%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.93%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.93%"
